$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy style (bold border) from A16 to new rows A17:A19
$ws.Range("A16").Copy($ws.Range("A17:A19"))

# Row 10: A=8 "Gaussian-Quadrature"
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = "Gaussian-Quadrature"
$ws.Cells.Item(10, 3).Value = 1.210448842380348
$ws.Cells.Item(10, 4).Value = 0.9085960483277356
$ws.Cells.Item(10, 5).Value = 1.05656414440336
$ws.Cells.Item(10, 6).Value = 0.9104558492700081
$ws.Cells.Item(10, 7).Value = 1.210448842380348
$ws.Cells.Item(10, 8).Value = 0.9085960483277356
$ws.Cells.Item(10, 9).Value = 1.067917046248743
$ws.Cells.Item(10, 10).Value = 0.91967140095761
$ws.Cells.Item(10, 11).Value = 1.021511881579576
$ws.Cells.Item(10, 12).Value = 0.8813047786622972
$ws.Cells.Item(10, 13).Value = 1.210448842380348
$ws.Cells.Item(10, 14).Value = 0.9825800963655478
$ws.Cells.Item(10, 15).Value = 1.021516221095363
$ws.Cells.Item(10, 16).Value = 0.9970587489787096

# Row 11: A=9 "Spiral-90deg-10rot-5space"
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = "Spiral-90deg-10rot-5space"
$ws.Cells.Item(11, 3).Value = 1.040159271627852
$ws.Cells.Item(11, 4).Value = 0.2836683827155736
$ws.Cells.Item(11, 5).Value = 1.389797186059282
$ws.Cells.Item(11, 6).Value = 0.8739336620014455
$ws.Cells.Item(11, 7).Value = 1.040159271627852
$ws.Cells.Item(11, 8).Value = 0.2836683827155736
$ws.Cells.Item(11, 9).Value = 1.259594408488008
$ws.Cells.Item(11, 10).Value = 1.000859648104836
$ws.Cells.Item(11, 11).Value = 1.048142785377809
$ws.Cells.Item(11, 12).Value = 0.5774774865891594
$ws.Cells.Item(11, 13).Value = 1.040159271627852
$ws.Cells.Item(11, 14).Value = 0.8367327843874279
$ws.Cells.Item(11, 15).Value = 0.8968896256010386
$ws.Cells.Item(11, 16).Value = 0.9342041038704959

# Row 12: A=10 "Spiral-90deg-15rot-5space"
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = "Spiral-90deg-15rot-5space"
$ws.Cells.Item(12, 3).Value = 1.036540111628894
$ws.Cells.Item(12, 4).Value = 0.2843757162888689
$ws.Cells.Item(12, 5).Value = 1.38972404633226
$ws.Cells.Item(12, 6).Value = 0.8749847640722022
$ws.Cells.Item(12, 7).Value = 1.036540111628894
$ws.Cells.Item(12, 8).Value = 0.2843757162888689
$ws.Cells.Item(12, 9).Value = 1.258884866183629
$ws.Cells.Item(12, 10).Value = 1.001990075303949
$ws.Cells.Item(12, 11).Value = 1.047691576373158
$ws.Cells.Item(12, 12).Value = 0.5784904433601387
$ws.Cells.Item(12, 13).Value = 1.036540111628894
$ws.Cells.Item(12, 14).Value = 0.8370498813105642
$ws.Cells.Item(12, 15).Value = 0.8964061595805561
$ws.Cells.Item(12, 16).Value = 0.9340851999428874

# Row 13: A=11 "Spiral-90deg-10rot-3space"
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = "Spiral-90deg-10rot-3space"
$ws.Cells.Item(13, 3).Value = 1.03910062401136
$ws.Cells.Item(13, 4).Value = 0.283769873122901
$ws.Cells.Item(13, 5).Value = 1.390101781067999
$ws.Cells.Item(13, 6).Value = 0.8741251775600941
$ws.Cells.Item(13, 7).Value = 1.03910062401136
$ws.Cells.Item(13, 8).Value = 0.283769873122901
$ws.Cells.Item(13, 9).Value = 1.259507847376276
$ws.Cells.Item(13, 10).Value = 1.0012076086197
$ws.Cells.Item(13, 11).Value = 1.04784550087967
$ws.Cells.Item(13, 12).Value = 0.5775629411149235
$ws.Cells.Item(13, 13).Value = 1.03910062401136
$ws.Cells.Item(13, 14).Value = 0.8369358270954499
$ws.Cells.Item(13, 15).Value = 0.8967743639405883
$ws.Cells.Item(13, 16).Value = 0.9341526692191153

# Row 14: A=12 "NoRotation-tilt60deg"
$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = "NoRotation-tilt60deg"
$ws.Cells.Item(14, 3).Value = 0.809355999999999
$ws.Cells.Item(14, 4).Value = 0.2099199999999999
$ws.Cells.Item(14, 5).Value = 1.574979999999996
$ws.Cells.Item(14, 6).Value = 0.8733199999999987
$ws.Cells.Item(14, 7).Value = 0.809355999999999
$ws.Cells.Item(14, 8).Value = 0.2099199999999999
$ws.Cells.Item(14, 9).Value = 1.344795999999999
$ws.Cells.Item(14, 10).Value = 0.9952480000000019
$ws.Cells.Item(14, 11).Value = 1.074148
$ws.Cells.Item(14, 12).Value = 0.5459639999999998
$ws.Cells.Item(14, 13).Value = 0.809355999999999
$ws.Cells.Item(14, 14).Value = 0.8924499999999979
$ws.Cells.Item(14, 15).Value = 0.8668939999999983
$ws.Cells.Item(14, 16).Value = 0.9284664999999992

# Row 15: A=13 "Rotation-NoTilt"
$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = "Rotation-NoTilt"
$ws.Cells.Item(15, 3).Value = 1.019162500000001
$ws.Cells.Item(15, 4).Value = 0.08
$ws.Cells.Item(15, 5).Value = 1.62
$ws.Cells.Item(15, 6).Value = 0.8025125000000032
$ws.Cells.Item(15, 7).Value = 1.019162500000001
$ws.Cells.Item(15, 8).Value = 0.08
$ws.Cells.Item(15, 9).Value = 1.444024999999998
$ws.Cells.Item(15, 10).Value = 0.867924999999999
$ws.Cells.Item(15, 11).Value = 1.180687500000001
$ws.Cells.Item(15, 12).Value = 0.4180750000000004
$ws.Cells.Item(15, 13).Value = 1.019162500000001
$ws.Cells.Item(15, 14).Value = 0.8500000000000001
$ws.Cells.Item(15, 15).Value = 0.880418750000001
$ws.Cells.Item(15, 16).Value = 0.9290484375000003

# Row 16: A=14 "Rotation-60detTilt"
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "Rotation-60detTilt"
$ws.Cells.Item(16, 3).Value = 1.0182387278848
$ws.Cells.Item(16, 4).Value = 0.4578196030464002
$ws.Cells.Item(16, 5).Value = 1.353303950131197
$ws.Cells.Item(16, 6).Value = 0.8822623808512001
$ws.Cells.Item(16, 7).Value = 1.0182387278848
$ws.Cells.Item(16, 8).Value = 0.4578196030464002
$ws.Cells.Item(16, 9).Value = 1.2507053677568
$ws.Cells.Item(16, 10).Value = 0.9215856194560031
$ws.Cells.Item(16, 11).Value = 1.099477903871996
$ws.Cells.Item(16, 12).Value = 0.6589874265088005
$ws.Cells.Item(16, 13).Value = 1.018224414822399
$ws.Cells.Item(16, 14).Value = 0.9055617765887987
$ws.Cells.Item(16, 15).Value = 0.9279061654783993
$ws.Cells.Item(16, 16).Value = 0.9552976224383997

# Row 17: A=15 "HexGrid-90degTilt5degRes"
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "HexGrid-90degTilt5degRes"
$ws.Cells.Item(17, 3).Value = 0.9954741158333702
$ws.Cells.Item(17, 4).Value = 0.9956309456830806
$ws.Cells.Item(17, 5).Value = 0.9916696509680184
$ws.Cells.Item(17, 6).Value = 0.99096221622741
$ws.Cells.Item(17, 7).Value = 0.9954741158333702
$ws.Cells.Item(17, 8).Value = 0.9956309456830806
$ws.Cells.Item(17, 9).Value = 0.992325672973139
$ws.Cells.Item(17, 10).Value = 0.9932343334042283
$ws.Cells.Item(17, 11).Value = 0.9926970675175544
$ws.Cells.Item(17, 12).Value = 0.991181935518539
$ws.Cells.Item(17, 13).Value = 0.9954665401410732
$ws.Cells.Item(17, 14).Value = 0.9936502983255495
$ws.Cells.Item(17, 15).Value = 0.9934342321779698
$ws.Cells.Item(17, 16).Value = 0.9928969922656674

# Row 18: A=16 "HexGrid-90degTilt22p5degRes"
$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 2).Value = "HexGrid-90degTilt22p5degRes"
$ws.Cells.Item(18, 3).Value = 0.9555999066360479
$ws.Cells.Item(18, 4).Value = 0.9920333835649312
$ws.Cells.Item(18, 5).Value = 0.9897523578422132
$ws.Cells.Item(18, 6).Value = 1.014038933559427
$ws.Cells.Item(18, 7).Value = 0.9555999066360479
$ws.Cells.Item(18, 8).Value = 0.9920333835649312
$ws.Cells.Item(18, 9).Value = 0.9774125866171012
$ws.Cells.Item(18, 10).Value = 1.005624824689247
$ws.Cells.Item(18, 11).Value = 0.9952870229557815
$ws.Cells.Item(18, 12).Value = 1.006358908091816
$ws.Cells.Item(18, 13).Value = 0.9555999066360479
$ws.Cells.Item(18, 14).Value = 0.9908928707035722
$ws.Cells.Item(18, 15).Value = 0.9878561454006549
$ws.Cells.Item(18, 16).Value = 0.9920134904945705

# Row 19: A=17 "HexGrid-60degTilt5degRes"
$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(19, 2).Value = "HexGrid-60degTilt5degRes"
$ws.Cells.Item(19, 3).Value = 0.9435945384509613
$ws.Cells.Item(19, 4).Value = 1.086309105141223
$ws.Cells.Item(19, 5).Value = 0.9727803913028665
$ws.Cells.Item(19, 6).Value = 1.00840885990544
$ws.Cells.Item(19, 7).Value = 0.9435945384509613
$ws.Cells.Item(19, 8).Value = 1.086309105141223
$ws.Cells.Item(19, 9).Value = 0.9558379423795529
$ws.Cells.Item(19, 10).Value = 1.00848941600661
$ws.Cells.Item(19, 11).Value = 0.9683419697192863
$ws.Cells.Item(19, 12).Value = 1.053131090150523
$ws.Cells.Item(19, 13).Value = 0.9435561134809031
$ws.Cells.Item(19, 14).Value = 1.029544748222045
$ws.Cells.Item(19, 15).Value = 1.002773223700123
$ws.Cells.Item(19, 16).Value = 0.999611664132058
